$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 134
$ws.Range("I33").Value = 124.28571
$ws.Range("K33").Value = 124.28571
$ws.Range("M33").Value = 104.71429
$ws.Range("H40").Value = 2635.8
$ws.Range("I40").Value = 2373.111
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2373.111
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2198.111
$ws.Range("N40").Value = -5350
$ws.Range("I76").Value = 3500
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3500
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3185
$ws.Range("N76").ClearContents()
$ws.Range("I79").Value = 3500
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3500
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2408
$ws.Range("N79").ClearContents()
$ws.Range("H98").Value = 2868.8215
$ws.Range("I98").Value = 1727.8636
$ws.Range("K98").Value = 1727.8636
$ws.Range("M98").Value = -229.8635999999999
$ws.Range("H106").Value = 2156.5
$ws.Range("I106").Value = 1880.5555
$ws.Range("K106").Value = 1880.5555
$ws.Range("M106").Value = -1249.5555
$ws.Range("H122").Value = 2868.8215
$ws.Range("I122").Value = 1727.8636
$ws.Range("K122").Value = 5183.5908
$ws.Range("M122").Value = -2733.5908
$ws.Range("H129").Value = 1946.0769
$ws.Range("I129").Value = 1471.8572
$ws.Range("J129").Value = 2499.3333
$ws.Range("K129").Value = 4415.571599999999
$ws.Range("L129").Value = 7497.999899999999
$ws.Range("M129").Value = 584.4284000000007
$ws.Range("N129").Value = -17497.9999
$ws.Range("H132").Value = 7938884.5
$ws.Range("I132").Value = 9525895
$ws.Range("K132").Value = 28577685
$ws.Range("M132").Value = -28575155
$ws.Range("H137").Value = 19486.723
$ws.Range("I137").Value = 2046.25
$ws.Range("J137").Value = 24469.715
$ws.Range("K137").Value = 6138.75
$ws.Range("L137").Value = 73409.145
$ws.Range("M137").Value = -3588.75
$ws.Range("N137").Value = -78509.145
$ws.Range("H138").Value = 1323604.2
$ws.Range("I138").Value = 2335.2778
$ws.Range("J138").Value = 1903673.5
$ws.Range("K138").Value = 7005.8334
$ws.Range("L138").Value = 5711020.5
$ws.Range("M138").Value = -1865.8334
$ws.Range("N138").Value = -5721300.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34763.758
$ws.Range("I32").Value = 36444.516
$ws.Range("K32").Value = 36444.516
$ws.Range("M32").Value = -36157.516
$ws.Range("H61").Value = 7431.9355
$ws.Range("I61").Value = 4606.364
$ws.Range("K61").Value = 4606.364
$ws.Range("M61").Value = -4394.364
$ws.Range("H63").Value = 4228.706
$ws.Range("J63").Value = 5400
$ws.Range("L63").Value = 5400
$ws.Range("N63").Value = -6772
$ws.Range("H66").Value = 4228.706
$ws.Range("J66").Value = 5400
$ws.Range("L66").Value = 27000
$ws.Range("N66").Value = -33864
$ws.Range("H110").Value = 20387.818
$ws.Range("I110").Value = 23582.81
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 23582.81
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = -21537.81
$ws.Range("N110").Value = -7590
$ws.Range("H122").Value = 1743.6316
$ws.Range("I122").Value = 1324.625
$ws.Range("K122").Value = 3973.875
$ws.Range("M122").Value = -1523.875
$ws.Range("H136").Value = 7431.9355
$ws.Range("I136").Value = 4606.364
$ws.Range("K136").Value = 13819.092
$ws.Range("M136").Value = -11269.092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3198.75
$ws.Range("I99").Value = 2098.3333
$ws.Range("K99").Value = 2098.3333
$ws.Range("M99").Value = -600.3332999999998
$ws.Range("H105").Value = 2428.85
$ws.Range("I105").Value = 2167.375
$ws.Range("K105").Value = 2167.375
$ws.Range("M105").Value = -420.375
$ws.Range("H107").Value = 915.4167
$ws.Range("I107").Value = 781.2353000000001
$ws.Range("J107").Value = 1241.2858
$ws.Range("K107").Value = 781.2353000000001
$ws.Range("L107").Value = 1241.2858
$ws.Range("M107").Value = 1138.7647
$ws.Range("N107").Value = -5081.2858
$ws.Range("H126").Value = 110000
$ws.Range("J126").Value = 110000
$ws.Range("L126").Value = 110000
$ws.Range("N126").Value = -119880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1025
$ws.Range("J22").Value = 1750
$ws.Range("L22").Value = 1750
$ws.Range("N22").Value = -2450
$ws.Range("H31").Value = 1820411
$ws.Range("I31").Value = 3848337
$ws.Range("J31").Value = 2270.6553
$ws.Range("K31").Value = 3848337
$ws.Range("L31").Value = 2270.6553
$ws.Range("M31").Value = -3848042
$ws.Range("N31").Value = -2860.6553
$ws.Range("H34").Value = 1820411
$ws.Range("I34").Value = 3848337
$ws.Range("J34").Value = 2270.6553
$ws.Range("K34").Value = 3848337
$ws.Range("L34").Value = 2270.6553
$ws.Range("M34").Value = -3848135
$ws.Range("N34").Value = -2674.6553
$ws.Range("H58").Value = 1255.4166
$ws.Range("I58").Value = 823.92
$ws.Range("K58").Value = 823.92
$ws.Range("M58").Value = -620.92
$ws.Range("H99").Value = 7371.3125
$ws.Range("I99").Value = 7195.7
$ws.Range("K99").Value = 7195.7
$ws.Range("M99").Value = -5697.7
$ws.Range("H105").Value = 1245.4615
$ws.Range("I105").Value = 622
$ws.Range("K105").Value = 622
$ws.Range("M105").Value = 1125
$ws.Range("H126").Value = 7371.3125
$ws.Range("I126").Value = 7195.7
$ws.Range("K126").Value = 21587.1
$ws.Range("M126").Value = -19117.1
$ws.Range("H132").Value = 2386.25
$ws.Range("I132").Value = 1814.4375
$ws.Range("J132").Value = 4673.5
$ws.Range("K132").Value = 5443.3125
$ws.Range("L132").Value = 14020.5
$ws.Range("M132").Value = -2913.3125
$ws.Range("N132").Value = -19080.5
$ws.Range("H136").Value = 1255.4166
$ws.Range("I136").Value = 823.92
$ws.Range("K136").Value = 2471.76
$ws.Range("M136").Value = 78.24000000000024

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 685.2381
$ws.Range("J122").Value = 729.94116
$ws.Range("L122").Value = 6569.47044
$ws.Range("N122").Value = -11469.47044
$ws.Range("H129").Value = 3315.4546
$ws.Range("J129").Value = 4055.8572
$ws.Range("L129").Value = 12167.5716
$ws.Range("N129").Value = -22167.5716
$ws.Range("H131").Value = 2925.3928
$ws.Range("I131").Value = 4733
$ws.Range("K131").Value = 14199
$ws.Range("M131").Value = -9159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1818.4762
$ws.Range("I102").Value = 1640.7368
$ws.Range("K102").Value = 1640.7368
$ws.Range("M102").Value = -18.7367999999999
$ws.Range("H113").Value = 1126.909
$ws.Range("I113").Value = 889.6
$ws.Range("K113").Value = 889.6
$ws.Range("M113").Value = 1280.4
$ws.Range("H130").Value = 162499.5
$ws.Range("J130").Value = 162499.5
$ws.Range("L130").Value = 162499.5
$ws.Range("N130").Value = -172539.5
$ws.Range("H132").Value = 1808.0233
$ws.Range("I132").Value = 1570.4103
$ws.Range("K132").Value = 4711.2309
$ws.Range("M132").Value = -2181.2309

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1411.25
$ws.Range("J16").Value = 1474
$ws.Range("L16").Value = 1474
$ws.Range("N16").Value = -1814
$ws.Range("H46").Value = 3823.5417
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 4188.25
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 4188.25
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -4564.25
$ws.Range("H55").Value = 2126.9412
$ws.Range("J55").Value = 2784.889
$ws.Range("L55").Value = 2784.889
$ws.Range("N55").Value = -3130.889
$ws.Range("H122").Value = 3411.8572
$ws.Range("I122").Value = 3045.75
$ws.Range("K122").Value = 9137.25
$ws.Range("M122").Value = -6687.25
$ws.Range("H134").Value = 127979.5
$ws.Range("J134").Value = 127979.5
$ws.Range("L134").Value = 127979.5
$ws.Range("N134").Value = -138119.5
$ws.Range("H136").Value = 5252.4
$ws.Range("I136").Value = 2981.5
$ws.Range("J136").Value = 6766.3335
$ws.Range("K136").Value = 8944.5
$ws.Range("L136").Value = 20299.0005
$ws.Range("M136").Value = -6394.5
$ws.Range("N136").Value = -25399.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 622.587
$ws.Range("I113").Value = 672.1111
$ws.Range("J113").Value = 444.3
$ws.Range("K113").Value = 2016.3333
$ws.Range("L113").Value = 1332.9
$ws.Range("M113").Value = 153.6667000000002
$ws.Range("N113").Value = -5672.9
$ws.Range("H122").Value = 5281.6787
$ws.Range("I122").Value = 6085.8423
$ws.Range("J122").Value = 3584
$ws.Range("K122").Value = 18257.5269
$ws.Range("L122").Value = 10752
$ws.Range("M122").Value = -15807.5269
$ws.Range("N122").Value = -15652
$ws.Range("H126").Value = 3577.111
$ws.Range("I126").Value = 2844
$ws.Range("K126").Value = 8532
$ws.Range("M126").Value = -6062
$ws.Range("H132").Value = 19979.717
$ws.Range("I132").Value = 21505.88
$ws.Range("J132").Value = 3955
$ws.Range("K132").Value = 64517.64
$ws.Range("L132").Value = 11865
$ws.Range("M132").Value = -61987.64
$ws.Range("N132").Value = -16925
$ws.Range("H133").Value = 115000
$ws.Range("J133").Value = 115000
$ws.Range("L133").Value = 115000
$ws.Range("N133").Value = -125120
